$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sugerir 2 Features (Francisco)"
$ws.Range("A3").Value = "Sugerir 2 Features (Iago)"
$ws.Range("A4").Value = "Sugerir 2 Features (James)"
$ws.Range("A5").Value = "Sugerir 2 Features (Joao)"
$ws.Range("A6").Value = "Sugerir 2 Features (Ricardo)"

$ws.Range("A2:A6").Select()
